$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = '26.951.67'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '1.809.36'
$ws.Range("E3").Value = '  -0.60%  '
Set-TextValue $ws.Range("D4") '1.001'
$ws.Range("E4").Value = '  -0.07%  '
Set-TextValue $ws.Range("D5") '310.64'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("E6").Value = '  -0.02%  '
Set-TextValue $ws.Range("D7") '0.4640'
$ws.Range("E7").Value = '  +3.97%  '
Set-TextValue $ws.Range("D8") '0.3719'
$ws.Range("E8").Value = '  -1.41%  '
Set-TextValue $ws.Range("D9") '0.07377'
$ws.Range("E9").Value = '  -0.26%  '
Set-TextValue $ws.Range("D10") '0.8759'
$ws.Range("E10").Value = '  -0.30%  '
Set-TextValue $ws.Range("D11") '20.47'
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("D12").Value = '1.804.46'
$ws.Range("E12").Value = '  -0.86%  '
Set-TextValue $ws.Range("D13") '5.375'
$ws.Range("E13").Value = '  -0.73%  '
Set-TextValue $ws.Range("D14") '92.54'
$ws.Range("E14").Value = '  -0.58%  '
Set-TextValue $ws.Range("D15") '6.530'
$ws.Range("E15").Value = '  -2.68%  '
Set-TextValue $ws.Range("D16") '0.07061'
$ws.Range("E16").Value = '  -0.79%  '
Set-TextValue $ws.Range("D17") '1.002'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("E18").Value = '  -0.85%  '
Set-TextValue $ws.Range("D19") '1.001'
$ws.Range("E19").Value = '  -0.04%  '
Set-TextValue $ws.Range("D20") '14.73'
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("D21").Value = '26.943.19'
$ws.Range("E21").Value = '  -1.07%  '
Set-TextValue $ws.Range("D22") '5.309'
$ws.Range("E22").Value = '  -0.83%  '
Set-TextValue $ws.Range("D23") '10.64'
$ws.Range("E23").Value = '  -2.43%  '
$ws.Range("D24").Value = '2.045.43'
$ws.Range("E24").Value = '  -0.24%  '
Set-TextValue $ws.Range("D25") '1.911'
$ws.Range("E25").Value = '  -2.70%  '
Set-TextValue $ws.Range("D26") '151.67'
$ws.Range("E26").Value = '  +0.41%  '
Set-TextValue $ws.Range("D27") '18.41'
$ws.Range("E27").Value = '  -0.84%  '
Set-TextValue $ws.Range("D28") '2.152'
$ws.Range("E28").Value = '  -5.99%  '
Set-TextValue $ws.Range("D29") '5.295'
$ws.Range("E29").Value = '  -0.82%  '
Set-TextValue $ws.Range("D30") '116.01'
$ws.Range("E30").Value = '  -1.11%  '
Set-TextValue $ws.Range("D31") '0.08929'
$ws.Range("E31").Value = '  +0.72%  '
Set-TextValue $ws.Range("D32") '0.7573'
$ws.Range("E32").Value = '  -3.19%  '
$ws.Range("E33").Value = '  -2.84%  '
Set-TextValue $ws.Range("D34") '2.934'
$ws.Range("E34").Value = '  +0.90%  '
Set-TextValue $ws.Range("D35") '4.461'
$ws.Range("E35").Value = '  -2.27%  '
Set-TextValue $ws.Range("D36") '1.001'
$ws.Range("E36").Value = '  -0.05%  '
Set-TextValue $ws.Range("D37") '1.106'
$ws.Range("E37").Value = '  -0.19%  '
Set-TextValue $ws.Range("D38") '0.01971'
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  -0.12%  '
Set-TextValue $ws.Range("D40") '2.429'
$ws.Range("E40").Value = '  +7.06%  '
Set-TextValue $ws.Range("D41") '2.924'
$ws.Range("E41").Value = '  +2.18%  '
Set-TextValue $ws.Range("D42") '7.221'
$ws.Range("E42").Value = '  -1.27%  '
Set-TextValue $ws.Range("D43") '0.5307'
$ws.Range("E43").Value = '  +0.33%  '
Set-TextValue $ws.Range("D44") '0.1667'
$ws.Range("E44").Value = '  -2.12%  '
Set-TextValue $ws.Range("D45") '8.522'
$ws.Range("E45").Value = '  -0.73%  '
Set-TextValue $ws.Range("D46") '0.5004'
$ws.Range("E46").Value = '  -0.44%  '
Set-TextValue $ws.Range("D47") '10.41'
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D48") '1.677'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D49") '104.05'
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D50") '1.001'
$ws.Range("E50").Value = '  +0.02%  '
Set-TextValue $ws.Range("D51") '0.06297'
$ws.Range("E51").Value = '  -1.40%  '
